$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change ("Generate Report for handoff") drops the second source file
# (a5482964-...) from the report - it had no localized content yet - and
# refreshes the handoff status/timestamps for the remaining source file
# (6378e7a4-...).
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1) Update status text + handoff timestamps on the row that stays (row 2).
#    "Handed back: in sync with en-US" is the very same shared string used by
#    the Overview rollup and by the per-locale "Status" column, so all of its
#    occurrences turn into "Ready for handoff" together.
$wsOverview.Cells.Item(2, 2).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(2, 3).Value2 = "Ready for handoff"
$wsZhCn.Cells.Item(2, 2).Value2 = "Ready for handoff"
$wsDeDe.Cells.Item(2, 2).Value2 = "Ready for handoff"

$wsZhCn.Cells.Item(2, 4).Value2 = "2016-01-11 13:14:03"
$wsDeDe.Cells.Item(2, 4).Value2 = "2016-01-11 13:14:21"

# 2) Remove the a5482964-... row (row 3) from every sheet - it shifts the
#    ".localization-config" row up from row 4 to row 3.
$wsOverview.Rows.Item(3).Delete()
$wsZhCn.Rows.Item(3).Delete()
$wsDeDe.Rows.Item(3).Delete()

# 3) The engine does not re-target hyperlinks when rows shift, so rebuild the
#    hyperlink collections on each sheet from scratch to match the new rows.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/.localization-config", "", "", ".localization-config") | Out-Null

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77aaee05489d57a998f90c9fbeeccda263acc5c8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a453083726a0c160f2e6a7b4812a50f36b58638c/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a852c14680e4aa1fe68b2441ef97771f3fb8265a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/.localization-config", "", "", ".localization-config") | Out-Null

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b4eb38c350b3a17ebca0689a29cb3ec1148dc86/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/857ea7a8730fbf661938fc764cc33b9960c8f889/e2e/6378e7a4-4c34-4a45-987c-e3baeb12303f.md", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5b04dabd619dc65744e05253c36c71275ee4dd7f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf", "", "", "6378e7a4-4c34-4a45-987c-e3baeb12303f.721f826828e83cc1026f124e830d1456e79c0502.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed82e45f74970dc7f5cfd688aaab708b214388f/.localization-config", "", "", ".localization-config") | Out-Null
